# Insert a new weekly record above row 59, shifting the existing
# rows 59:121 down to 60:122, then populate the new row 59 with the
# new observation (same variety/unit/origin/volume as the former
# row 59, but an updated date and updated prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values that must carry over unchanged into the new row
# before the insert shifts everything down.
$h59 = $ws.Range("H59").Value()
$j59 = $ws.Range("J59").Value()
$n59 = $ws.Range("N59").Value()
$o59 = $ws.Range("O59").Value()
$q59 = $ws.Range("Q59").Value()

# Shift rows 59:121 down to 60:122, inserting a new blank row 59.
$ws.Rows("59:59").Insert()

# Populate the new row 59 (columns A-C, E, F, G, I, R match every
# other row in this sheet, so copy them from the row right below).
$ws.Range("A59").Value = $ws.Range("A60").Value()
$ws.Range("B59").Value = $ws.Range("B60").Value()
$ws.Range("C59").Value = $ws.Range("C60").Value()
$ws.Range("D59").Value = 44930
$ws.Range("E59").Value = $ws.Range("E60").Value()
$ws.Range("F59").Value = $ws.Range("F60").Value()
$ws.Range("G59").Value = $ws.Range("G60").Value()
$ws.Range("H59").Value = $h59
$ws.Range("I59").Value = $ws.Range("I60").Value()
$ws.Range("J59").Value = $j59
$ws.Range("K59").Value = 13000
$ws.Range("L59").Value = 14000
$ws.Range("M59").Value = 13500
$ws.Range("N59").Value = $n59
$ws.Range("O59").Value = $o59
$ws.Range("P59").Value = 900
$ws.Range("Q59").Value = $q59
$ws.Range("R59").Value = $ws.Range("R60").Value()
